$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- New column D: "CO2 (2017) in g/kWh End " ----------------------------
# Header F1 first (copy format from A4 so it picks up the same "apply"
# flags as the rest of the table, then re-style it to a bold label with no
# border), so that the new cell-style records land in the same creation
# order as the authored workbook.
$ws.Range("A4").Copy() | Out-Null
$ws.Range("F1").PasteSpecial(-4122) | Out-Null
$ws.Range("F1").Value = "Quelle 2017"
$ws.Range("F1").Font.Name = "MS Sans Serif"
$ws.Range("F1").Font.Size = 10
$ws.Range("F1").Font.Bold = $true
$ws.Range("F1").Borders.LineStyle = 0

$ws.Range("D1").Value = "CO2 (2017) in g/kWh End "
$ws.Range("D1").Font.Bold = $true
$ws.Range("D1").Borders.LineStyle = 1
$ws.Range("D1").Borders.Weight = 2
$ws.Range("D1").Borders.Color = 0

$ws.Range("D2").Value = 310
$ws.Range("D2").Font.Name = "Calibri"
$ws.Range("D2").Font.Size = 11
$ws.Range("D2").Font.Bold = $false
$ws.Range("D2").Borders.LineStyle = 1
$ws.Range("D2").Borders.Weight = 2
$ws.Range("D2").Borders.Color = 0

$ws.Range("D3").Value = 242
$ws.Range("D4").Value = 258
$ws.Range("D5").Value = 419
$ws.Range("D6").Value = 417
$ws.Range("D7").Value = 12
$ws.Range("D8").Value = 12
$ws.Range("D9").Value = 12
$ws.Range("D10").Value = 494
$ws.Range("D11").Value = 62
$ws.Range("D12").Value = 840

$ws.Range("D3:D12").Font.Name = "Calibri"
$ws.Range("D3:D12").Font.Size = 11
$ws.Range("D3:D12").Font.Bold = $false
$ws.Range("D3:D12").Borders.LineStyle = 1
$ws.Range("D3:D12").Borders.Weight = 2
$ws.Range("D3:D12").Borders.Color = 0

# --- Source note in F2 -----------------------------------------------------
$ws.Range("F2").Value = "https://www.ifeu.de/energie/pdf/ifeu_Endbericht_Weiterentwicklung_PEF.pdf"

# --- Column width for the new column D -------------------------------------
$ws.Columns.Item(4).ColumnWidth = 23.28515625

# --- Page setup --------------------------------------------------------------
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# --- Selection ---------------------------------------------------------------
$ws.Range("F2").Select() | Out-Null
